$d = $word.ActiveDocument

$replacements = @(
    @{old = "345÷8=43, 1"; new = "640÷8=80, 0"},
    @{old = "224÷8=28, 0"; new = "656÷9=72, 8"},
    @{old = "873÷3=291, 0"; new = "771÷2=385, 1"},
    @{old = "611÷9=67, 8"; new = "996÷2=498, 0"},
    @{old = "639÷9=71, 0"; new = "352÷2=176, 0"},
    @{old = "563÷6=93, 5"; new = "128÷8=16, 0"},
    @{old = "200÷5=40, 0"; new = "113÷3=37, 2"},
    @{old = "365÷2=182, 1"; new = "506÷3=168, 2"},
    @{old = "547÷8=68, 3"; new = "537÷4=134, 1"},
    @{old = "865÷5=173, 0"; new = "290÷4=72, 2"},
    @{old = "377÷4=94, 1"; new = "455÷6=75, 5"},
    @{old = "716÷4=179, 0"; new = "816÷5=163, 1"},
    @{old = "740÷6=123, 2"; new = "363÷2=181, 1"},
    @{old = "711÷6=118, 3"; new = "304÷5=60, 4"},
    @{old = "363÷6=60, 3"; new = "967÷3=322, 1"},
    @{old = "657÷2=328, 1"; new = "864÷4=216, 0"},
    @{old = "995÷5=199, 0"; new = "143÷9=15, 8"},
    @{old = "560÷8=70, 0"; new = "613÷6=102, 1"},
    @{old = "768÷7=109, 5"; new = "537÷3=179, 0"},
    @{old = "239÷5=47, 4"; new = "695÷9=77, 2"},
    @{old = "972÷4=243, 0"; new = "342÷9=38, 0"},
    @{old = "932÷2=466, 0"; new = "119÷9=13, 2"},
    @{old = "479÷7=68, 3"; new = "434÷8=54, 2"},
    @{old = "965÷8=120, 5"; new = "645÷5=129, 0"},
    @{old = "568÷7=81, 1"; new = "695÷3=231, 2"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
